$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Priority"

$ws.Range("A14").Value = "Classification"
$ws.Range("B14").Value = "Choose the 'private' option"
$ws.Range("C14").Value = 'Your .ics will be created with a "private" visibility setting.'

$ws.Range("C14").Select() | Out-Null
